# Bugs Fixed and Cooling Effect Done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 8.380000000000001
$ws.Range("E11").Value = 8.44
$ws.Range("E12").Value = 8.51
$ws.Range("E13").Value = 7.32
$ws.Range("E14").Value = 4.62
$ws.Range("E32").Value = 38.09999999999999
